$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.35%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'35.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.92%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.108"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.95%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08083"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.03%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.948"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.36%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'4.211"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.84%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'7.751"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.03%"
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'0.90%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1389"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'13.77%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.1912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.54%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.09193"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.46%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.03460"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.59%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.09833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.26%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001412"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.59%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.005774"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.47%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'3.619"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'3.40%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.15%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.3440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.84%"
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.1341"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.64%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'4.914"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-2.49%"
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'-0.88%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.04435"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.45%"
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.001222"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.59%"
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.004829"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.47%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0001242"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.63%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D39").Value = "'0.02013"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.30%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.04922"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.66%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.007718"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.26%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.01013"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.19%"
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'3.49%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.002104"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.30%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.01156"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.55%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00006466"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.05%"
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.18%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D49").Value = "'0.001193"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-8.50%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.18%"
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.0002004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.18%"
$ws.Range("E51").Style = "Normal"
